# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" between "2021-Q1" and "总计", holding
#   the per-fund holdings detail for the new quarter.
# - Insert a new row at the top of "总计"'s data for the 2022-Q1 summary,
#   pushing the existing 2021-Q1 summary row down to row 3.

$wb = $excel.ActiveWorkbook

$q1_2021 = $wb.Worksheets.Item(1)   # "2021-Q1" stays sheet 1, untouched

# ---------------------------------------------------------------------
# 1. New "2022-Q1" worksheet, positioned right after "2021-Q1"
# ---------------------------------------------------------------------
$q1_2022 = $wb.Worksheets.Add($null, $q1_2021)
$q1_2022.Name = "2022-Q1"

# NOTE: worksheet references obtained before Worksheets.Add() can read back
# stale/blank values afterwards, so (re)fetch every sheet handle we need
# *after* the Add() call.
$q1_2022 = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")

# Pull header formatting (bold font + border) and the index-column format
# from the "总计" sheet so the new sheet matches the workbook's look.
$totalSheet.Range("B1").Copy()
$q1_2022.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$totalSheet.Range("A2").Copy()
$q1_2022.Range("A2:A3").PasteSpecial(-4122)   # xlPasteFormats

# Write the General-formatted cells (numeric index column + rank column)
# before the text columns get NumberFormat="@" applied below, otherwise
# they lose the formatting that was just pasted onto them.
$q1_2022.Range("A2").Value = 0
$q1_2022.Range("A3").Value = 1
$q1_2022.Range("H1").Value = "仓位排名"
$q1_2022.Range("H2").Value = 9
$q1_2022.Range("H3").Value = 9

# Store the numeric-looking text columns (fund code, size, weight, ...) as
# plain text so leading zeros / exact decimal strings are preserved.
$q1_2022.Range("B1:G3").NumberFormat = "@"

$q1_2022.Range("B1").Value = "基金代码"
$q1_2022.Range("C1").Value = "基金名称"
$q1_2022.Range("D1").Value = "基金规模"
$q1_2022.Range("E1").Value = "股票总仓位"
$q1_2022.Range("F1").Value = "仓位占比"
$q1_2022.Range("G1").Value = "持有市值(亿元)"

$q1_2022.Range("B2").Value = "009613"
$q1_2022.Range("C2").Value = "上银中证500指数增强A"
$q1_2022.Range("D2").Value = "2.83"
$q1_2022.Range("E2").Value = "90.41"
$q1_2022.Range("F2").Value = "1.20"
$q1_2022.Range("G2").Value = "0.0340"

$q1_2022.Range("B3").Value = "009614"
$q1_2022.Range("C3").Value = "上银中证500指数增强C"
$q1_2022.Range("D3").Value = "1.70"
$q1_2022.Range("E3").Value = "90.41"
$q1_2022.Range("F3").Value = "1.20"
$q1_2022.Range("G3").Value = "0.0204"

# ---------------------------------------------------------------------
# 2. "总计" sheet: add the 2022-Q1 summary row above the 2021-Q1 one
# ---------------------------------------------------------------------
# Move the existing 2021-Q1 row down to row 3 (literal values, so the
# new row 2 can be written fresh without disturbing formatting). Grab
# row 2's current index-column style (still the original "2021-Q1" A2
# formatting at this point) before it gets overwritten below.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q1"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.14

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.05

Write-Host "2022-Q1 sheet added and summary sheet updated"
